$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.988.19"
$ws.Cells.Item(2, 5).Value = "  +0.19%  "
$ws.Cells.Item(3, 4).Value = "'2.639.23"
$ws.Cells.Item(3, 5).Value = "  +0.35%  "
$ws.Cells.Item(4, 5).Value = "  -0.05%  "
$ws.Cells.Item(5, 4).Value = "'579.80"
$ws.Cells.Item(5, 5).Value = "  +0.34%  "
$ws.Cells.Item(6, 4).Value = "'156.77"
$ws.Cells.Item(6, 5).Value = "  +0.72%  "
$ws.Cells.Item(7, 5).Value = "  -3.25%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 4).Value = "'2.637.60"
$ws.Cells.Item(9, 5).Value = "  +0.31%  "
$ws.Cells.Item(10, 5).Value = "  -2.63%  "
$ws.Cells.Item(12, 5).Value = "  -0.85%  "
$ws.Cells.Item(13, 5).Value = "  +0.84%  "
$ws.Cells.Item(14, 4).Value = "'28.72"
$ws.Cells.Item(14, 5).Value = "  +0.46%  "
$ws.Cells.Item(15, 4).Value = "'3.117.31"
$ws.Cells.Item(15, 5).Value = "  +0.18%  "
$ws.Cells.Item(16, 4).Value = "'0.0000185"
$ws.Cells.Item(16, 5).Value = "  +0.08%  "
$ws.Cells.Item(17, 4).Value = "'63.898.17"
$ws.Cells.Item(17, 5).Value = "  +0.20%  "
$ws.Cells.Item(18, 4).Value = "'2.650.57"
$ws.Cells.Item(18, 5).Value = "  -0.19%  "
$ws.Cells.Item(19, 4).Value = "'12.15"
$ws.Cells.Item(19, 5).Value = "  -0.21%  "
$ws.Cells.Item(20, 4).Value = "'7.77"
$ws.Cells.Item(20, 5).Value = "  +2.24%  "
$ws.Cells.Item(21, 5).Value = "  -2.92%  "
$ws.Cells.Item(22, 4).Value = "'345.18"
$ws.Cells.Item(22, 5).Value = "  -0.45%  "
$ws.Cells.Item(23, 5).Value = "  +0.19%  "
$ws.Cells.Item(24, 5).Value = "  +0.81%  "
$ws.Cells.Item(25, 5).Value = "  +7.59%  "
$ws.Cells.Item(26, 4).Value = "'0.0000112"
$ws.Cells.Item(26, 5).Value = "  +2.88%  "
$ws.Cells.Item(27, 4).Value = "'9.30"
$ws.Cells.Item(27, 5).Value = "  -0.45%  "
$ws.Cells.Item(28, 4).Value = "'1.63"
$ws.Cells.Item(28, 5).Value = "  +3.52%  "
$ws.Cells.Item(29, 4).Value = "'581.59"
$ws.Cells.Item(29, 5).Value = "  +1.30%  "
$ws.Cells.Item(30, 4).Value = "'8.19"
$ws.Cells.Item(30, 5).Value = "  +3.27%  "
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 4).Value = "'0.999"
$ws.Cells.Item(33, 5).Value = "  -1.14%  "
$ws.Cells.Item(34, 4).Value = "'1.73"
$ws.Cells.Item(34, 5).Value = "  +1.14%  "
$ws.Cells.Item(35, 4).Value = "'6.62"
$ws.Cells.Item(35, 5).Value = "  +2.35%  "
$ws.Cells.Item(36, 4).Value = "'5.46"
$ws.Cells.Item(36, 5).Value = "  +2.84%  "
$ws.Cells.Item(37, 5).Value = "  -2.05%  "
$ws.Cells.Item(38, 4).Value = "'19.78"
$ws.Cells.Item(38, 5).Value = "  -1.01%  "
$ws.Cells.Item(39, 5).Value = "  -0.06%  "
$ws.Cells.Item(40, 5).Value = "  +1.87%  "
$ws.Cells.Item(41, 4).Value = "'153.23"
$ws.Cells.Item(41, 5).Value = "  +0.59%  "
$ws.Cells.Item(42, 4).Value = "'2.55"
$ws.Cells.Item(42, 5).Value = "  +7.65%  "
$ws.Cells.Item(43, 5).Value = "  -0.01%  "
$ws.Cells.Item(44, 4).Value = "'162.41"
$ws.Cells.Item(44, 5).Value = "  +2.09%  "
$ws.Cells.Item(45, 4).Value = "'24.16"
$ws.Cells.Item(45, 5).Value = "  +4.30%  "
$ws.Cells.Item(46, 4).Value = "'3.91"
$ws.Cells.Item(46, 5).Value = "  -1.89%  "
$ws.Cells.Item(47, 4).Value = "'0.0589"
$ws.Cells.Item(47, 5).Value = "  -1.46%  "
$ws.Cells.Item(48, 5).Value = "  +0.37%  "
$ws.Cells.Item(49, 5).Value = "  -2.22%  "
$ws.Cells.Item(50, 5).Value = "  -1.59%  "
$ws.Cells.Item(51, 4).Value = "'19.11"
$ws.Cells.Item(51, 5).Value = "  -0.20%  "
